$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '40.080.24'
$ws.Range("E2").Value = '  +0.11%  '

# Row 3
$ws.Range("D3").Value = '2.220.57'
$ws.Range("E3").Value = '  -0.87%  '

# Row 4
$ws.Range("E4").Value = '  +0.12%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '291.54'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.71%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '86.98'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.21%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.511'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.62%  '

# Row 8
$ws.Range("E8").Value = '  +0.12%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.470'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.85%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '30.16'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.14%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0778'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.02%  '

# Row 12
$ws.Range("E12").Value = '  +2.99%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.45'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.04%  '

# Row 14
$ws.Range("D14").Value = '2.568.31'
$ws.Range("E14").Value = '  -0.68%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '13.83'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.95%  '

# Row 16
$ws.Range("D16").Value = '2.234.88'
$ws.Range("E16").Value = '  +0.15%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.725'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.08%  '

# Row 18
$ws.Range("D18").Value = '40.000.75'
$ws.Range("E18").Value = '  +0.12%  '

# Row 19
$ws.Range("D19").Value = '0.0₃0883'
$ws.Range("E19").Value = '  -1.33%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.26'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +6.23%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.78'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.23%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '65.49'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.41%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '235.84'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.29%  '

# Row 24
$ws.Range("E24").Value = '  -0.11%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.44'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.77%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.80'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.11%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.65'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.04%  '

# Row 28
$ws.Range("E28").Value = '  -1.79%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.18'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.46%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '156.06'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.76%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '31.51'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -8.09%  '

# Row 32
$ws.Range("E32").Value = '  -0.05%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.91'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.14%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0715'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.03%  '

# Row 35
$ws.Range("B35").Value = 'WEMIXToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.35'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.38%  '

# Row 36
$ws.Range("B36").Value = 'LidoDAOToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.88'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +5.91%  '

# Row 37
$ws.Range("E37").Value = '  +0.33%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '15.59'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -5.87%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0975'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.33%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.68'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.02%  '

# Row 41
$ws.Range("D41").Value = '2.121.36'
$ws.Range("E41").Value = '  +7.04%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.83'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.83%  '

# Row 43
$ws.Range("B43").Value = 'ApeXProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.14'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -4.04%  '

# Row 44
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '18.16'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +10.82%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0267'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.85%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '9.79'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.54%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.63'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.02%  '

# Row 48
$ws.Range("D48").Value = '2.435.61'
$ws.Range("E48").Value = '  -0.85%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.45'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.03%  '

# Row 50
$ws.Range("B50").Value = 'TrustWalletToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.10'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.82%  '

# Row 51
$ws.Range("B51").Value = 'BitcoinSV'
$ws.Range("C51").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '68.75'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.50%  '
